# Adds the "24. 8. 2021" wave column (AH on "data", AG on "pocetR")
# and bumps the "aktualizace" footer date from 27. 7. 2021 to 1. 9. 2021.

$wb = $excel.ActiveWorkbook

# --- Sheet "data": new column AH (34), header + 110 data rows ---
$ws1 = $wb.Worksheets.Item("data")

$ws1.Cells.Item(1, 34).Value = "24. 8. 2021"
$ws1.Range("AG1").Copy()
$ws1.Range("AH1").PasteSpecial(-4122)

$ws1.Cells.Item(2, 34).Value = 0.16
$ws1.Cells.Item(3, 34).Value = 0.16
$ws1.Cells.Item(4, 34).Value = 0.22
$ws1.Cells.Item(5, 34).Value = 0.15
$ws1.Cells.Item(6, 34).Value = 0.31
$ws1.Cells.Item(7, 34).Value = 0.15
$ws1.Cells.Item(8, 34).Value = 0.07000000000000001
$ws1.Cells.Item(9, 34).Value = 0.23
$ws1.Cells.Item(10, 34).Value = 0.14
$ws1.Cells.Item(11, 34).Value = 0.41
$ws1.Cells.Item(12, 34).Value = 0.11
$ws1.Cells.Item(13, 34).Value = 0.13
$ws1.Cells.Item(14, 34).Value = 0.21
$ws1.Cells.Item(15, 34).Value = 0.17
$ws1.Cells.Item(16, 34).Value = 0.38
$ws1.Cells.Item(17, 34).Value = 0.21
$ws1.Cells.Item(18, 34).Value = 0.23
$ws1.Cells.Item(19, 34).Value = 0.21
$ws1.Cells.Item(20, 34).Value = 0.15
$ws1.Cells.Item(21, 34).Value = 0.2
$ws1.Cells.Item(22, 34).Value = 0.1
$ws1.Cells.Item(23, 34).Value = 0.11
$ws1.Cells.Item(24, 34).Value = 0.27
$ws1.Cells.Item(25, 34).Value = 0.11
$ws1.Cells.Item(26, 34).Value = 0.41
$ws1.Cells.Item(27, 34).Value = 0.19
$ws1.Cells.Item(28, 34).Value = 0.14
$ws1.Cells.Item(29, 34).Value = 0.22
$ws1.Cells.Item(30, 34).Value = 0.13
$ws1.Cells.Item(31, 34).Value = 0.32
$ws1.Cells.Item(32, 34).Value = 0.18
$ws1.Cells.Item(33, 34).Value = 0.24
$ws1.Cells.Item(34, 34).Value = 0.16
$ws1.Cells.Item(35, 34).Value = 0.22
$ws1.Cells.Item(36, 34).Value = 0.2
$ws1.Cells.Item(37, 34).Value = 0.07000000000000001
$ws1.Cells.Item(38, 34).Value = 0.13
$ws1.Cells.Item(39, 34).Value = 0.23
$ws1.Cells.Item(40, 34).Value = 0.15
$ws1.Cells.Item(41, 34).Value = 0.42
$ws1.Cells.Item(42, 34).Value = 0.14
$ws1.Cells.Item(43, 34).Value = 0.19
$ws1.Cells.Item(44, 34).Value = 0.19
$ws1.Cells.Item(45, 34).Value = 0.19
$ws1.Cells.Item(46, 34).Value = 0.29
$ws1.Cells.Item(47, 34).Value = 0.19
$ws1.Cells.Item(48, 34).Value = 0.15
$ws1.Cells.Item(49, 34).Value = 0.22
$ws1.Cells.Item(50, 34).Value = 0.14
$ws1.Cells.Item(51, 34).Value = 0.3
$ws1.Cells.Item(52, 34).Value = 0.15
$ws1.Cells.Item(53, 34).Value = 0.12
$ws1.Cells.Item(54, 34).Value = 0.21
$ws1.Cells.Item(55, 34).Value = 0.17
$ws1.Cells.Item(56, 34).Value = 0.35
$ws1.Cells.Item(57, 34).Value = 0.17
$ws1.Cells.Item(58, 34).Value = 0.19
$ws1.Cells.Item(59, 34).Value = 0.22
$ws1.Cells.Item(60, 34).Value = 0.13
$ws1.Cells.Item(61, 34).Value = 0.29
$ws1.Cells.Item(62, 34).Value = 0.15
$ws1.Cells.Item(63, 34).Value = 0.15
$ws1.Cells.Item(64, 34).Value = 0.26
$ws1.Cells.Item(65, 34).Value = 0.12
$ws1.Cells.Item(66, 34).Value = 0.32
$ws1.Cells.Item(67, 34).Value = 0.17
$ws1.Cells.Item(68, 34).Value = 0.16
$ws1.Cells.Item(69, 34).Value = 0.17
$ws1.Cells.Item(70, 34).Value = 0.17
$ws1.Cells.Item(71, 34).Value = 0.33
$ws1.Cells.Item(72, 34).Value = 0.14
$ws1.Cells.Item(73, 34).Value = 0.18
$ws1.Cells.Item(74, 34).Value = 0.17
$ws1.Cells.Item(75, 34).Value = 0.22
$ws1.Cells.Item(76, 34).Value = 0.29
$ws1.Cells.Item(77, 34).Value = 0.19
$ws1.Cells.Item(78, 34).Value = 0.17
$ws1.Cells.Item(79, 34).Value = 0.18
$ws1.Cells.Item(80, 34).Value = 0.19
$ws1.Cells.Item(81, 34).Value = 0.27
$ws1.Cells.Item(82, 34).Value = 0.18
$ws1.Cells.Item(83, 34).Value = 0.03
$ws1.Cells.Item(84, 34).Value = 0.25
$ws1.Cells.Item(85, 34).Value = 0.16
$ws1.Cells.Item(86, 34).Value = 0.38
$ws1.Cells.Item(87, 34).Value = 0.1
$ws1.Cells.Item(88, 34).Value = 0.08
$ws1.Cells.Item(89, 34).Value = 0.2
$ws1.Cells.Item(90, 34).Value = 0.19
$ws1.Cells.Item(91, 34).Value = 0.43
$ws1.Cells.Item(92, 34).Value = 0.18
$ws1.Cells.Item(93, 34).Value = 0.23
$ws1.Cells.Item(94, 34).Value = 0.2
$ws1.Cells.Item(95, 34).Value = 0.16
$ws1.Cells.Item(96, 34).Value = 0.23
$ws1.Cells.Item(97, 34).Value = 0.12
$ws1.Cells.Item(98, 34).Value = 0.1
$ws1.Cells.Item(99, 34).Value = 0.21
$ws1.Cells.Item(100, 34).Value = 0.11
$ws1.Cells.Item(101, 34).Value = 0.46
$ws1.Cells.Item(102, 34).Value = 0.12
$ws1.Cells.Item(103, 34).Value = 0.18
$ws1.Cells.Item(104, 34).Value = 0.22
$ws1.Cells.Item(105, 34).Value = 0.13
$ws1.Cells.Item(106, 34).Value = 0.35
$ws1.Cells.Item(107, 34).Value = 0.24
$ws1.Cells.Item(108, 34).Value = 0.24
$ws1.Cells.Item(109, 34).Value = 0.22
$ws1.Cells.Item(110, 34).Value = 0.14
$ws1.Cells.Item(111, 34).Value = 0.16

$ws1.Range("A112").Value = "Život během pandemie, Kontakty vs. protektivní aktivity, % respondentů celkově a ve skupinách, aktualizace 1. 9. 2021"

# --- Sheet "pocetR": new column AG (33), header + 22 data rows ---
$ws2 = $wb.Worksheets.Item("pocetR")

$ws2.Cells.Item(1, 33).Value = "24. 8. 2021"
$ws2.Range("AF1").Copy()
$ws2.Range("AG1").PasteSpecial(-4122)

$ws2.Cells.Item(2, 33).Value = 1638
$ws2.Cells.Item(3, 33).Value = 374
$ws2.Cells.Item(4, 33).Value = 584
$ws2.Cells.Item(5, 33).Value = 680
$ws2.Cells.Item(6, 33).Value = 463
$ws2.Cells.Item(7, 33).Value = 708
$ws2.Cells.Item(8, 33).Value = 464
$ws2.Cells.Item(9, 33).Value = 267
$ws2.Cells.Item(10, 33).Value = 316
$ws2.Cells.Item(11, 33).Value = 1055
$ws2.Cells.Item(12, 33).Value = 827
$ws2.Cells.Item(13, 33).Value = 811
$ws2.Cells.Item(14, 33).Value = 839
$ws2.Cells.Item(15, 33).Value = 388
$ws2.Cells.Item(16, 33).Value = 193
$ws2.Cells.Item(17, 33).Value = 218
$ws2.Cells.Item(18, 33).Value = 199
$ws2.Cells.Item(19, 33).Value = 323
$ws2.Cells.Item(20, 33).Value = 306
$ws2.Cells.Item(21, 33).Value = 175
$ws2.Cells.Item(22, 33).Value = 262
$ws2.Cells.Item(23, 33).Value = 373

$ws2.Range("A24").Value = "Život během pandemie, Kontakty vs. protektivní aktivity, velikost dotázaného souboru celkově a ve skupinách, aktualizace 1. 9. 2021"
$ws2.Cells.Item(24, 33).Value = "'"
$ws2.Cells.Item(24, 33).ClearFormats()
